$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.274280333333333
$ws.Range("H2").Value = 6.822841
$ws.Range("I2").Value = 0.2367408419877493
$ws.Range("J2").Value = 0.2367408419877492
$ws.Range("M2").Value = 0.7339303333333334
$ws.Range("N2").Value = 2.201791
$ws.Range("O2").Value = 0.03574007706012852
$ws.Range("P2").Value = 0.03574007706012852
$ws.Range("Q2").Value = 1.669163323136778
$ws.Range("R2").Value = 15.022469908231
$ws.Range("S2").Value = 0.008461135935921868
$ws.Range("T2").Value = 0.008461135935921866
$ws.Range("G3").Value = 2.274280333333333
$ws.Range("H3").Value = 6.822841
$ws.Range("I3").Value = 0.2367408419877493
$ws.Range("J3").Value = 0.2367408419877492
$ws.Range("O3").Value = 0.3842514532634088
$ws.Range("P3").Value = 0.3842514532634088
$ws.Range("Q3").Value = 17.94563653486933
$ws.Range("R3").Value = 161.510728813824
$ws.Range("S3").Value = 0.09096801258059568
$ws.Range("T3").Value = 0.09096801258059566
$ws.Range("G4").Value = 2.274280333333333
$ws.Range("H4").Value = 6.822841
$ws.Range("I4").Value = 0.2367408419877493
$ws.Range("J4").Value = 0.2367408419877492
$ws.Range("M4").Value = 4.974008666666667
$ws.Range("N4").Value = 14.922026
$ws.Range("O4").Value = 0.2422184299659874
$ws.Range("P4").Value = 0.2422184299659874
$ws.Range("Q4").Value = 11.31229008842956
$ws.Range("R4").Value = 101.810610795866
$ws.Range("S4").Value = 0.05734299505509854
$ws.Range("T4").Value = 0.05734299505509852
$ws.Range("G5").Value = 2.274280333333333
$ws.Range("H5").Value = 6.822841
$ws.Range("I5").Value = 0.2367408419877493
$ws.Range("J5").Value = 0.2367408419877492
$ws.Range("M5").Value = 2.087648
$ws.Range("N5").Value = 6.262943999999999
$ws.Range("O5").Value = 0.1016618294757629
$ws.Range("P5").Value = 0.1016618294757629
$ws.Range("Q5").Value = 4.747896789322667
$ws.Range("R5").Value = 42.731071103904
$ws.Range("S5").Value = 0.0240675071081071
$ws.Range("T5").Value = 0.02406750710810709
$ws.Range("G6").Value = 2.274280333333333
$ws.Range("H6").Value = 6.822841
$ws.Range("I6").Value = 0.2367408419877493
$ws.Range("J6").Value = 0.2367408419877492
$ws.Range("M6").Value = 4.848944666666667
$ws.Range("N6").Value = 14.546834
$ws.Range("O6").Value = 0.2361282102347124
$ws.Range("P6").Value = 0.2361282102347124
$ws.Range("Q6").Value = 11.02785949282156
$ws.Range("R6").Value = 99.25073543539401
$ws.Range("S6").Value = 0.05590119130802608
$ws.Range("T6").Value = 0.05590119130802607
$ws.Range("I7").Value = 0.3055660932301819
$ws.Range("J7").Value = 0.3055660932301818
$ws.Range("M7").Value = 0.7339303333333334
$ws.Range("N7").Value = 2.201791
$ws.Range("O7").Value = 0.03574007706012852
$ws.Range("P7").Value = 0.03574007706012852
$ws.Range("Q7").Value = 2.154422157712889
$ws.Range("R7").Value = 19.389799419416
$ws.Range("S7").Value = 0.01092095571900912
$ws.Range("T7").Value = 0.01092095571900911
$ws.Range("I8").Value = 0.3055660932301819
$ws.Range("J8").Value = 0.3055660932301818
$ws.Range("O8").Value = 0.3842514532634088
$ws.Range("P8").Value = 0.3842514532634088
$ws.Range("S8").Value = 0.1174142153917196
$ws.Range("T8").Value = 0.1174142153917196
$ws.Range("I9").Value = 0.3055660932301819
$ws.Range("J9").Value = 0.3055660932301818
$ws.Range("M9").Value = 4.974008666666667
$ws.Range("N9").Value = 14.922026
$ws.Range("O9").Value = 0.2422184299659874
$ws.Range("P9").Value = 0.2422184299659874
$ws.Range("Q9").Value = 14.60099684864178
$ws.Range("R9").Value = 131.408971637776
$ws.Range("S9").Value = 0.07401373935305519
$ws.Range("T9").Value = 0.07401373935305516
$ws.Range("I10").Value = 0.3055660932301819
$ws.Range("J10").Value = 0.3055660932301818
$ws.Range("M10").Value = 2.087648
$ws.Range("N10").Value = 6.262943999999999
$ws.Range("O10").Value = 0.1016618294757629
$ws.Range("P10").Value = 0.1016618294757629
$ws.Range("Q10").Value = 6.128204414549333
$ws.Range("R10").Value = 55.15383973094399
$ws.Range("S10").Value = 0.03106440806354182
$ws.Range("T10").Value = 0.03106440806354181
$ws.Range("I11").Value = 0.3055660932301819
$ws.Range("J11").Value = 0.3055660932301818
$ws.Range("M11").Value = 4.848944666666667
$ws.Range("N11").Value = 14.546834
$ws.Range("O11").Value = 0.2361282102347124
$ws.Range("P11").Value = 0.2361282102347124
$ws.Range("Q11").Value = 14.23387664595378
$ws.Range("R11").Value = 128.104889813584
$ws.Range("S11").Value = 0.07215277470285611
$ws.Range("T11").Value = 0.07215277470285608
$ws.Range("G12").Value = 2.075403
$ws.Range("H12").Value = 6.226209000000001
$ws.Range("I12").Value = 0.2160387382692492
$ws.Range("J12").Value = 0.2160387382692492
$ws.Range("M12").Value = 0.7339303333333334
$ws.Range("N12").Value = 2.201791
$ws.Range("O12").Value = 0.03574007706012852
$ws.Range("P12").Value = 0.03574007706012852
$ws.Range("Q12").Value = 1.523201215591
$ws.Range("R12").Value = 13.708810940319
$ws.Range("S12").Value = 0.007721241153715901
$ws.Range("T12").Value = 0.0077212411537159
$ws.Range("G13").Value = 2.075403
$ws.Range("H13").Value = 6.226209000000001
$ws.Range("I13").Value = 0.2160387382692492
$ws.Range("J13").Value = 0.2160387382692492
$ws.Range("O13").Value = 0.3842514532634088
$ws.Range("P13").Value = 0.3842514532634088
$ws.Range("Q13").Value = 16.376357547264
$ws.Range("R13").Value = 147.387217925376
$ws.Range("S13").Value = 0.0830131991411522
$ws.Range("T13").Value = 0.0830131991411522
$ws.Range("G14").Value = 2.075403
$ws.Range("H14").Value = 6.226209000000001
$ws.Range("I14").Value = 0.2160387382692492
$ws.Range("J14").Value = 0.2160387382692492
$ws.Range("M14").Value = 4.974008666666667
$ws.Range("N14").Value = 14.922026
$ws.Range("O14").Value = 0.2422184299659874
$ws.Range("P14").Value = 0.2422184299659874
$ws.Range("Q14").Value = 10.323072508826
$ws.Range("R14").Value = 92.90765257943401
$ws.Range("S14").Value = 0.05232856399541042
$ws.Range("T14").Value = 0.05232856399541041
$ws.Range("G15").Value = 2.075403
$ws.Range("H15").Value = 6.226209000000001
$ws.Range("I15").Value = 0.2160387382692492
$ws.Range("J15").Value = 0.2160387382692492
$ws.Range("M15").Value = 2.087648
$ws.Range("N15").Value = 6.262943999999999
$ws.Range("O15").Value = 0.1016618294757629
$ws.Range("P15").Value = 0.1016618294757629
$ws.Range("Q15").Value = 4.332710922144
$ws.Range("R15").Value = 38.994398299296
$ws.Range("S15").Value = 0.02196289337008739
$ws.Range("T15").Value = 0.02196289337008738
$ws.Range("G16").Value = 2.075403
$ws.Range("H16").Value = 6.226209000000001
$ws.Range("I16").Value = 0.2160387382692492
$ws.Range("J16").Value = 0.2160387382692492
$ws.Range("M16").Value = 4.848944666666667
$ws.Range("N16").Value = 14.546834
$ws.Range("O16").Value = 0.2361282102347124
$ws.Range("P16").Value = 0.2361282102347124
$ws.Range("Q16").Value = 10.063514308034
$ws.Range("R16").Value = 90.57162877230601
$ws.Range("S16").Value = 0.05101284060888327
$ws.Range("T16").Value = 0.05101284060888327
$ws.Range("G17").Value = 0.9024383333333333
$ws.Range("H17").Value = 2.707315
$ws.Range("I17").Value = 0.09393917176526073
$ws.Range("J17").Value = 0.09393917176526072
$ws.Range("M17").Value = 0.7339303333333334
$ws.Range("N17").Value = 2.201791
$ws.Range("O17").Value = 0.03574007706012852
$ws.Range("P17").Value = 0.03574007706012852
$ws.Range("Q17").Value = 0.6623268667961111
$ws.Range("R17").Value = 5.960941801165
$ws.Range("S17").Value = 0.003357393237855068
$ws.Range("T17").Value = 0.003357393237855067
$ws.Range("G18").Value = 0.9024383333333333
$ws.Range("H18").Value = 2.707315
$ws.Range("I18").Value = 0.09393917176526073
$ws.Range("J18").Value = 0.09393917176526072
$ws.Range("O18").Value = 0.3842514532634088
$ws.Range("P18").Value = 0.3842514532634088
$ws.Range("Q18").Value = 7.120859327573333
$ws.Range("R18").Value = 64.08773394815999
$ws.Range("S18").Value = 0.03609626326916242
$ws.Range("T18").Value = 0.0360962632691624
$ws.Range("G19").Value = 0.9024383333333333
$ws.Range("H19").Value = 2.707315
$ws.Range("I19").Value = 0.09393917176526073
$ws.Range("J19").Value = 0.09393917176526072
$ws.Range("M19").Value = 4.974008666666667
$ws.Range("N19").Value = 14.922026
$ws.Range("O19").Value = 0.2422184299659874
$ws.Range("P19").Value = 0.2422184299659874
$ws.Range("Q19").Value = 4.488736091132223
$ws.Range("R19").Value = 40.39862482019
$ws.Range("S19").Value = 0.02275379869728667
$ws.Range("T19").Value = 0.02275379869728666
$ws.Range("G20").Value = 0.9024383333333333
$ws.Range("H20").Value = 2.707315
$ws.Range("I20").Value = 0.09393917176526073
$ws.Range("J20").Value = 0.09393917176526072
$ws.Range("M20").Value = 2.087648
$ws.Range("N20").Value = 6.262943999999999
$ws.Range("O20").Value = 0.1016618294757629
$ws.Range("P20").Value = 0.1016618294757629
$ws.Range("Q20").Value = 1.883973581706667
$ws.Range("R20").Value = 16.95576223536
$ws.Range("S20").Value = 0.00955002806109434
$ws.Range("T20").Value = 0.009550028061094337
$ws.Range("G21").Value = 0.9024383333333333
$ws.Range("H21").Value = 2.707315
$ws.Range("I21").Value = 0.09393917176526073
$ws.Range("J21").Value = 0.09393917176526072
$ws.Range("M21").Value = 4.848944666666667
$ws.Range("N21").Value = 14.546834
$ws.Range("O21").Value = 0.2361282102347124
$ws.Range("P21").Value = 0.2361282102347124
$ws.Range("Q21").Value = 4.375873543412223
$ws.Range("R21").Value = 39.38286189071
$ws.Range("S21").Value = 0.02218168849986224
$ws.Range("T21").Value = 0.02218168849986224
$ws.Range("G22").Value = 1.419044
$ws.Range("H22").Value = 4.257132
$ws.Range("I22").Value = 0.1477151547475591
$ws.Range("J22").Value = 0.1477151547475591
$ws.Range("M22").Value = 0.7339303333333334
$ws.Range("N22").Value = 2.201791
$ws.Range("O22").Value = 0.03574007706012852
$ws.Range("P22").Value = 0.03574007706012852
$ws.Range("Q22").Value = 1.041479435934667
$ws.Range("R22").Value = 9.373314923412002
$ws.Range("S22").Value = 0.005279351013626572
$ws.Range("T22").Value = 0.00527935101362657
$ws.Range("G23").Value = 1.419044
$ws.Range("H23").Value = 4.257132
$ws.Range("I23").Value = 0.1477151547475591
$ws.Range("J23").Value = 0.1477151547475591
$ws.Range("O23").Value = 0.3842514532634088
$ws.Range("P23").Value = 0.3842514532634088
$ws.Range("Q23").Value = 11.197233462272
$ws.Range("R23").Value = 100.775101160448
$ws.Range("S23").Value = 0.05675976288077891
$ws.Range("T23").Value = 0.05675976288077889
$ws.Range("G24").Value = 1.419044
$ws.Range("H24").Value = 4.257132
$ws.Range("I24").Value = 0.1477151547475591
$ws.Range("J24").Value = 0.1477151547475591
$ws.Range("M24").Value = 4.974008666666667
$ws.Range("N24").Value = 14.922026
$ws.Range("O24").Value = 0.2422184299659874
$ws.Range("P24").Value = 0.2422184299659874
$ws.Range("Q24").Value = 7.058337154381335
$ws.Range("R24").Value = 63.52503438943201
$ws.Range("S24").Value = 0.03577933286513664
$ws.Range("T24").Value = 0.03577933286513663
$ws.Range("G25").Value = 1.419044
$ws.Range("H25").Value = 4.257132
$ws.Range("I25").Value = 0.1477151547475591
$ws.Range("J25").Value = 0.1477151547475591
$ws.Range("M25").Value = 2.087648
$ws.Range("N25").Value = 6.262943999999999
$ws.Range("O25").Value = 0.1016618294757629
$ws.Range("P25").Value = 0.1016618294757629
$ws.Range("Q25").Value = 2.962464368512
$ws.Range("R25").Value = 26.662179316608
$ws.Range("S25").Value = 0.01501699287293229
$ws.Range("T25").Value = 0.01501699287293228
$ws.Range("G26").Value = 1.419044
$ws.Range("H26").Value = 4.257132
$ws.Range("I26").Value = 0.1477151547475591
$ws.Range("J26").Value = 0.1477151547475591
$ws.Range("M26").Value = 4.848944666666667
$ws.Range("N26").Value = 14.546834
$ws.Range("O26").Value = 0.2361282102347124
$ws.Range("P26").Value = 0.2361282102347124
$ws.Range("Q26").Value = 6.880865835565335
$ws.Range("R26").Value = 61.927792520088
$ws.Range("S26").Value = 0.03487971511508471
$ws.Range("T26").Value = 0.0348797151150847
